$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2539
$ws.Range("F7").Value = 1340
$ws.Range("F8").Value = 1700
$ws.Range("F9").Value = 183
$ws.Range("F11").Value = 2376
$ws.Range("F13").Value = 155
$ws.Range("F14").Value = 53
$ws.Range("F16").Value = 108
$ws.Range("F17").Value = 97
$ws.Range("F18").Value = 8615
$ws.Range("F20").Value = 6719
$ws.Range("F21").Value = 10849
$ws.Range("F24").Value = 211
$ws.Range("F26").Value = 528
$ws.Range("F27").Value = 756
$ws.Range("F28").Value = 198
$ws.Range("F29").Value = 174
$ws.Range("F30").Value = 2197
$ws.Range("F31").Value = 88
$ws.Range("F32").Value = 23
$ws.Range("F33").Value = 4447
$ws.Range("F34").Value = 453
$ws.Range("F35").Value = 439

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 13
$ws.Range("F8").Value = 1173

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 19

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 2539
$ws.Range("F10").Value = 1340
$ws.Range("F11").Value = 13
$ws.Range("F12").Value = 1700
$ws.Range("F14").Value = 183
$ws.Range("F15").Value = 2376
$ws.Range("F18").Value = 155
$ws.Range("F19").Value = 53
$ws.Range("F21").Value = 108
$ws.Range("F22").Value = 97
$ws.Range("F23").Value = 8615
$ws.Range("F25").Value = 6719
$ws.Range("F26").Value = 10850
$ws.Range("F30").Value = 211
$ws.Range("F32").Value = 528
$ws.Range("F36").Value = 198
$ws.Range("F37").Value = 174
$ws.Range("F38").Value = 23
$ws.Range("F39").Value = 4447
$ws.Range("F46").Value = 439
